# fix: corrige l'export excel de B3 avec l'ajout des 2 colonnes supplémentaires
#
# Sheet "B3" ("B3 - Énergie et émissions de gaz à effet de serre") is missing
# two columns in its "Consommation d'énergie par combustible" sub-table:
# "Densité" and "Valeur Calorifique Nette (NCV)". Insert them right before
# the existing "Quantité" column (column G) and label the new header cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("B3")

# Insert two new blank columns at G (pushes the former G.. columns to I..,
# and auto-extends/repositions the merged header ranges above them).
$ws.Columns.Item(7).Insert()
$ws.Columns.Item(7).Insert()

# Label the two new columns on the header row (row 4).
$ws.Cells.Item(4, 7).Value = "Densité"
$ws.Cells.Item(4, 8).Value = "Valeur Calorifique Nette (NCV)"

# Restore a sensible selection on the sheet (matches the author's saved view).
$ws.Range("J10").Select()
